$d = $word.ActiveDocument

$pairs = @(
    @("51×71=", "32×48="),
    @("49×84=", "84×63="),
    @("38×30=", "61×94="),
    @("56×55=", "24×82="),
    @("56×43=", "77×58="),
    @("24×76=", "22×94="),
    @("72×24=", "92×67="),
    @("53×68=", "30×36="),
    @("55×81=", "67×62="),
    @("71×43=", "16×42="),
    @("21×36=", "82×60="),
    @("87×70=", "58×96="),
    @("38×81=", "53×72="),
    @("70×57=", "22×77="),
    @("71×67=", "53×93="),
    @("40×69=", "93×28="),
    @("23×81=", "38×18="),
    @("71×77=", "64×81="),
    @("91×77=", "91×82="),
    @("68×14=", "25×87="),
    @("49×49=", "20×43="),
    @("47×87=", "36×33="),
    @("91×88=", "45×31="),
    @("98×73=", "74×72="),
    @("88×24=", "15×60=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
